# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to reflect the newly generated output (commit 456a3b4).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 381
$ws1.Range("F3").Value = 113
$ws1.Range("F4").Value = 1598
$ws1.Range("F5").Value = 17
$ws1.Range("F7").Value = 411
$ws1.Range("F8").Value = 142
$ws1.Range("F10").Value = 473

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 381
$ws4.Range("F3").Value = 113
$ws4.Range("F4").Value = 1598
$ws4.Range("F5").Value = 17
$ws4.Range("F6").Value = 23
$ws4.Range("F7").Value = 411
$ws4.Range("F8").Value = 142
$ws4.Range("F10").Value = 473
